$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11121.889
$ws.Range("I40").Value = 6148.5
$ws.Range("K40").Value = 6148.5
$ws.Range("M40").Value = -5973.5
$ws.Range("H111").Value = 828
$ws.Range("I111").Value = 450
$ws.Range("K111").Value = 1350
$ws.Range("M111").Value = 1717

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8374.727999999999
$ws.Range("I32").Value = 4545.732
$ws.Range("K32").Value = 4545.732
$ws.Range("M32").Value = -4258.732
$ws.Range("H45").Value = 3350
$ws.Range("I45").Value = 2775
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 2775
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -2398
$ws.Range("N45").Value = -5254
$ws.Range("H61").Value = 2325.611
$ws.Range("I61").Value = 1964.1333
$ws.Range("K61").Value = 1964.1333
$ws.Range("M61").Value = -1752.1333
$ws.Range("H74").Value = 3177.3333
$ws.Range("I74").Value = 2811
$ws.Range("J74").Value = 3273.7368
$ws.Range("K74").Value = 2811
$ws.Range("L74").Value = 3273.7368
$ws.Range("M74").Value = -1937
$ws.Range("N74").Value = -5021.736800000001
$ws.Range("H77").Value = 3177.3333
$ws.Range("I77").Value = 2811
$ws.Range("J77").Value = 3273.7368
$ws.Range("K77").Value = 14055
$ws.Range("L77").Value = 16368.684
$ws.Range("M77").Value = -9687
$ws.Range("N77").Value = -25104.684
$ws.Range("H122").Value = 4208.706
$ws.Range("I122").Value = 4475.385
$ws.Range("J122").Value = 3342
$ws.Range("K122").Value = 13426.155
$ws.Range("L122").Value = 10026
$ws.Range("M122").Value = -10976.155
$ws.Range("N122").Value = -14926
$ws.Range("H136").Value = 2325.611
$ws.Range("I136").Value = 1964.1333
$ws.Range("K136").Value = 5892.3999
$ws.Range("M136").Value = -3342.3999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 87430.914
$ws.Range("I22").Value = 125896.375
$ws.Range("J22").Value = 10500
$ws.Range("K22").Value = 125896.375
$ws.Range("L22").Value = 10500
$ws.Range("M22").Value = -125723.375
$ws.Range("N22").Value = -10846
$ws.Range("H132").Value = 81775
$ws.Range("J132").Value = 81775
$ws.Range("L132").Value = 81775
$ws.Range("N132").Value = -91895
$ws.Range("H134").Value = 1535.8536
$ws.Range("I134").Value = 1164.8572
$ws.Range("K134").Value = 3494.5716
$ws.Range("M134").Value = -959.5715999999998
$ws.Range("H138").Value = 96467.664
$ws.Range("J138").Value = 96467.664
$ws.Range("L138").Value = 96467.664
$ws.Range("N138").Value = -106747.664
$ws.Range("H140").Value = 99990
$ws.Range("J140").Value = 99990
$ws.Range("L140").Value = 99990
$ws.Range("N140").Value = -110350

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15095.038
$ws.Range("I31").Value = 1950.421
$ws.Range("K31").Value = 1950.421
$ws.Range("M31").Value = -1655.421
$ws.Range("H34").Value = 15095.038
$ws.Range("I34").Value = 1950.421
$ws.Range("K34").Value = 1950.421
$ws.Range("M34").Value = -1748.421
$ws.Range("H134").Value = 2062.3572
$ws.Range("I134").Value = 1763
$ws.Range("K134").Value = 5289
$ws.Range("M134").Value = -2754

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 63.25
$ws.Range("I15").Value = 60.5
$ws.Range("K15").Value = 181.5
$ws.Range("M15").Value = -41.5
$ws.Range("H39").Value = 7841.067
$ws.Range("I39").Value = 404.5
$ws.Range("J39").Value = 10545.272
$ws.Range("K39").Value = 1213.5
$ws.Range("L39").Value = 31635.816
$ws.Range("M39").Value = -919.5
$ws.Range("N39").Value = -32223.816
$ws.Range("H55").Value = 1400
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 4200
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -4023
$ws.Range("N55").ClearContents()
$ws.Range("H86").Value = 2785.1482
$ws.Range("J86").Value = 3175.3914
$ws.Range("L86").Value = 9526.174199999999
$ws.Range("N86").Value = -11898.1742
$ws.Range("H89").Value = 2785.1482
$ws.Range("J89").Value = 3175.3914
$ws.Range("L89").Value = 28578.5226
$ws.Range("N89").Value = -40434.5226
$ws.Range("H140").Value = 2914.923
$ws.Range("I140").Value = 1831.3334
$ws.Range("K140").Value = 5494.0002
$ws.Range("M140").Value = -314.0002000000004

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 282.46667
$ws.Range("I2").Value = 185
$ws.Range("J2").Value = 428.66666
$ws.Range("K2").Value = 185
$ws.Range("L2").Value = 428.66666
$ws.Range("M2").Value = -72
$ws.Range("N2").Value = -654.66666
$ws.Range("H70").Value = 234324.8
$ws.Range("I70").Value = 204684.12
$ws.Range("K70").Value = 204684.12
$ws.Range("M70").Value = -204414.12
$ws.Range("H73").Value = 234324.8
$ws.Range("I73").Value = 204684.12
$ws.Range("K73").Value = 204684.12
$ws.Range("M73").Value = -203748.12
$ws.Range("H80").Value = 2250
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 2250
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -22484
$ws.Range("H113").Value = 2779493
$ws.Range("I113").Value = 1111.8
$ws.Range("K113").Value = 1111.8
$ws.Range("M113").Value = 1058.2
$ws.Range("H126").Value = 3734.3845
$ws.Range("I126").Value = 2886.1667
$ws.Range("J126").Value = 4461.4287
$ws.Range("K126").Value = 8658.500100000001
$ws.Range("L126").Value = 13384.2861
$ws.Range("M126").Value = -6188.500100000001
$ws.Range("N126").Value = -18324.2861
$ws.Range("H132").Value = 4189.4062
$ws.Range("I132").Value = 3740.0833
$ws.Range("K132").Value = 11220.2499
$ws.Range("M132").Value = -8690.249899999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31830.5
$ws.Range("I7").Value = 15952.533
$ws.Range("K7").Value = 15952.533
$ws.Range("M7").Value = -15840.533
$ws.Range("H22").Value = 954.84375
$ws.Range("I22").Value = 849
$ws.Range("J22").Value = 1018.35
$ws.Range("K22").Value = 849
$ws.Range("L22").Value = 1018.35
$ws.Range("M22").Value = -554
$ws.Range("N22").Value = -1608.35
$ws.Range("H27").Value = 954.84375
$ws.Range("I27").Value = 849
$ws.Range("J27").Value = 1018.35
$ws.Range("K27").Value = 849
$ws.Range("L27").Value = 1018.35
$ws.Range("M27").Value = -742
$ws.Range("N27").Value = -1232.35
$ws.Range("H40").Value = 18521854
$ws.Range("I40").Value = 5002
$ws.Range("K40").Value = 5002
$ws.Range("M40").Value = -4866
$ws.Range("H100").Value = 37679.1
$ws.Range("I100").Value = 60149.668
$ws.Range("J100").Value = 3973.25
$ws.Range("K100").Value = 60149.668
$ws.Range("L100").Value = 3973.25
$ws.Range("M100").Value = -59608.668
$ws.Range("N100").Value = -5055.25
$ws.Range("H122").Value = 14320020
$ws.Range("I122").Value = 51142.445
$ws.Range("K122").Value = 153427.335
$ws.Range("M122").Value = -150977.335
$ws.Range("H126").Value = 31830.5
$ws.Range("I126").Value = 15952.533
$ws.Range("K126").Value = 47857.599
$ws.Range("M126").Value = -45387.599

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1787.7021
$ws.Range("I107").Value = 1263.0476
$ws.Range("K107").Value = 3789.142800000001
$ws.Range("M107").Value = -1869.142800000001
$ws.Range("H122").Value = 3195.0667
$ws.Range("I122").Value = 2710.375
$ws.Range("J122").Value = 3749
$ws.Range("K122").Value = 8131.125
$ws.Range("L122").Value = 11247
$ws.Range("M122").Value = -5681.125
$ws.Range("N122").Value = -16147
$ws.Range("H126").Value = 1326.1052
$ws.Range("I126").Value = 957.16
$ws.Range("J126").Value = 2035.6154
$ws.Range("K126").Value = 2871.48
$ws.Range("L126").Value = 6106.8462
$ws.Range("M126").Value = -401.48
$ws.Range("N126").Value = -11046.8462
$ws.Range("H136").Value = 2497.5557
$ws.Range("I136").Value = 2119.75
$ws.Range("J136").Value = 2799.8
$ws.Range("K136").Value = 6359.25
$ws.Range("L136").Value = 8399.400000000001
$ws.Range("M136").Value = -3809.25
$ws.Range("N136").Value = -13499.4
